# Edit: rewrite the "Baz chan" / "ges" paragraph (split across two runs
# around a _GoBack bookmark) into a single run reading
# "Hi Baz this is my pull request for the assignment 1", while keeping
# the bookmark intact and dropping the now-redundant trailing run.

$d = $word.ActiveDocument

# Locate the paragraph that still has the placeholder text, regardless
# of its position in the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Baz chan*") {
        $target = $p
    }
}

if ($target -ne $null) {
    # Replace just the first run's text ("Baz chan"); this leaves the
    # bookmarkStart/bookmarkEnd pair (and the trailing "ges" run) alone
    # because the match doesn't span them.
    $target.Range.Find.Execute("Baz chan", $true, $false, $false, $false, $false,
                                $true, 1, $false,
                                "Hi Baz this is my pull request for the assignment 1", 2)

    # Remove the now-orphaned trailing "ges" run by replacing it with
    # nothing, scoped to this paragraph only so it can't clobber other
    # occurrences of "ges" elsewhere in the document (e.g. "changes").
    $target.Range.Find.Execute("ges", $true, $false, $false, $false, $false,
                                $true, 1, $false, "", 2)
}
